# Add message and dialog xls files
# - Rename the existing "MessageReaderDB" ##var entry to "phone.MessageReaderDB"
# - Add a new ##var row describing the dialog reader db (mirrors the phone row above it)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# full_name for the phone message reader changes to include its module prefix
$ws.Range("B5").Value = "phone.MessageReaderDB"

# New row 6: dialog.DialogReaderDB / DiglogDB / read_schema_from_file=TRUE / dialog/mainDialogData.xlsx
$ws.Range("B6").Value = "dialog.DialogReaderDB"
$ws.Range("C6").Value = "DiglogDB"
$ws.Range("D6").Value = $true
$ws.Range("E6").Value = "dialog/mainDialogData.xlsx"

# Move the active selection to the newly added row, matching the post-edit cursor position
$null = $ws.Range("E6").Select()
